$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Gewichtung")
$ws2 = $wb.Worksheets.Item("Begründung")

# Remove the "total_weight" header and values (column G) - clear column G data rows 2-13
# but keep G1 styled (blank) - we just clear contents, formatting (style) is kept by Excel by default
$ws1.Range("G1").ClearContents()
$ws1.Range("G2:G13").ClearContents()

# New data values for Szenario 0
$data = @(
    @("umweltbelastung", "Vermeidung von Umweltbelastung", 0.25, "co2", "CO2-Emissionen", 0.5),
    @("umweltbelastung", "Vermeidung von Umweltbelastung", 0.25, "graueEnergie", "Graue Energie", 0.5),
    @("langlebigkeit_wirtschaftlichkeit", "Langlebigkeit und Wirtschaftlichkeit", 0.4, "lebensdauer", "Lebensdauer", 0.5),
    @("langlebigkeit_wirtschaftlichkeit", "Langlebigkeit und Wirtschaftlichkeit", 0.4, "unterhalt", "Unterhalt", 0.25),
    @("langlebigkeit_wirtschaftlichkeit", "Langlebigkeit und Wirtschaftlichkeit", 0.4, "kosten", "Kosten", 0.25),
    @("multifunktionale_nutzungsqualitaet", "Multifunktionale Nutzungsqualität", 0.1, "versickerung", "Versickerungsleistung", 0.25),
    @("multifunktionale_nutzungsqualitaet", "Multifunktionale Nutzungsqualität", 0.1, "oberflaechentemperatur", "Oberflächentemperatur", 0.25),
    @("multifunktionale_nutzungsqualitaet", "Multifunktionale Nutzungsqualität", 0.1, "befahrbarkeit", "Befahrbarkeit", 0.25),
    @("multifunktionale_nutzungsqualitaet", "Multifunktionale Nutzungsqualität", 0.1, "barrierefreiheit", "Barrierefreiheit", 0.25)
)

$r = 2
foreach ($row in $data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Rows 11-13: kreislauffaehigkeit group with formula 1/3
$ws1.Cells.Item(11, 1).Value = "kreislauffaehigkeit"
$ws1.Cells.Item(11, 2).Value = "Kreislauffähigkeit"
$ws1.Cells.Item(11, 3).Value = 0.25
$ws1.Cells.Item(11, 4).Value = "recyclingfaehigkeit"
$ws1.Cells.Item(11, 5).Value = "Recyclingfähigkeit"
$ws1.Range("F11").Formula = "=1/3"

$ws1.Cells.Item(12, 1).Value = "kreislauffaehigkeit"
$ws1.Cells.Item(12, 2).Value = "Kreislauffähigkeit"
$ws1.Cells.Item(12, 3).Value = 0.25
$ws1.Cells.Item(12, 4).Value = "lokaleMaterialien"
$ws1.Cells.Item(12, 5).Value = "Lokale Materialien"

$ws1.Cells.Item(13, 1).Value = "kreislauffaehigkeit"
$ws1.Cells.Item(13, 2).Value = "Kreislauffähigkeit"
$ws1.Cells.Item(13, 3).Value = 0.25
$ws1.Cells.Item(13, 4).Value = "modulareErneuerbarkeit"
$ws1.Cells.Item(13, 5).Value = "Modulare Erneuerbarkeit"

# Fill F12:F13 together as one operation so Excel stores it as a shared formula
$ws1.Range("F12:F13").Formula = "=1/3"

# Remove the "total_weight" shared string by deleting unused string (handled automatically by Excel when no longer referenced)

# Update sheet view selections / active sheet
# select Begründung's cell first, then Gewichtung last so Gewichtung ends up the active tab
$ws2.Range("C9").Select() | Out-Null
$ws1.Range("F25").Select() | Out-Null
